# Add a new "Introduction to data" course row before row 47 (becomes new row 46),
# pushing the existing rows 47-52 down to 48-53, and update the view/selection
# state to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 46; Excel shifts rows 46-52 down to 47-53 and the inserted
# row inherits formatting from the row above it (row 45), matching styles
# s="3" (A), s="102" (I) and s="2" (J/K) that the target file shows.
$null = $ws.Rows.Item(46).Insert()

# Populate the new row 46 with the "Introduction to data" course entry.
$ws.Cells.Item(46, 1).Value2 = "Introduction to data"   # A46 -> shared string "Introduction to data"
$ws.Cells.Item(46, 9).ClearContents()                    # I46 stays blank (keeps inherited style)
$ws.Cells.Item(46, 10).Clear()                            # J46 removed entirely (no value/style)
$ws.Cells.Item(46, 11).Value2 = 5                         # K46 = 5

# Update the view state: scrolled so row 33 is at the top, and the active
# selection sitting on the first empty row below the data (A54).
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
$null = $ws.Range("A54").Select()
